# Adds the "CNPJS", "ENDEREÇO BANCOS", "VALOR DE CAUSA" and "DATA" rows
# (18-21) to the reference sheet, formats the currency / date values,
# widens column A to fit the new labels and moves the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 18 - CNPJS
$ws.Range("A18").Value = "CNPJS"
$ws.Range("B18").Value = "33.010.284/0001-56, 19.041.377/0001-17, 90.261.648/0001-04, 22.611.216/0001-26, 45.279.812/0001-56, 13.902.317/0001-65"

# Row 19 - ENDEREÇO BANCOS
$ws.Range("A19").Value = "ENDEREÇO BANCOS"
$ws.Range("B19").Value = "R. Machado de Assis, Av. Bossa Nova, R. Canários, R. Catamarã, Av. Cristovão Colombo, Av. Oliveria Nunes, R. Parque das Águas"

# Row 20 - VALOR DE CAUSA (currency formatted number)
$ws.Range("A20").Value = "VALOR DE CAUSA"
$ws.Range("B20").Value = 30000
$ws.Range("B20").NumberFormat = '_-[$R$-416]\ * #,##0.00_-;\-[$R$-416]\ * #,##0.00_-;_-[$R$-416]\ * "-"??_-;_-@_-'

# Row 21 - DATA (date formatted value, 44776 serial == 2022-08-03)
$ws.Range("A21").Value = "DATA"
$ws.Range("B21").Value = Get-Date -Year 2022 -Month 8 -Day 3 -Hour 0 -Minute 0 -Second 0
$ws.Range("B21").NumberFormat = "mm-dd-yy"

# Widen column A so the new, longer labels fit (matches bestFit width of 16)
$ws.Columns.Item(1).ColumnWidth = 15.2

# Update the selected/active cell shown when the workbook is reopened
$ws.Range("E12").Select() | Out-Null
